$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4
$ws.Range("A4").Value = 112027322
$ws.Range("B4").Value = 90689
$ws.Range("D4").Value = "NT"
$ws.Range("E4").Value = 5966
$ws.Range("F4").Value = "Motaggsvamp"
$ws.Range("G4").Value = "Sarcodon squamosus"
$ws.Range("H4").Value = "(Schaeff.) Quél."
$ws.Range("Q4").Value = 752287.2680111516
$ws.Range("R4").Value = 7093027.007953994

# Row 5
$ws.Range("A5").Value = 112027339
$ws.Range("B5").Value = 90666
$ws.Range("D5").Value = "LC"
$ws.Range("E5").Value = 4364
$ws.Range("F5").Value = "Dropptaggsvamp"
$ws.Range("G5").Value = "Hydnellum ferrugineum"
$ws.Range("H5").Value = "(Fr.:Fr.) P. Karst."
$ws.Range("Q5").Value = 752300.1425363926
$ws.Range("R5").Value = 7092944.398308391
$ws.Range("AC5").Value = "Stora fruktkroppar"

# Row 6
$ws.Range("A6").Value = 112027371
$ws.Range("B6").Value = 90660
$ws.Range("D6").Value = "NT"
$ws.Range("E6").Value = 4362
$ws.Range("F6").Value = "Blå taggsvamp"
$ws.Range("G6").Value = "Hydnellum caeruleum"
$ws.Range("H6").Value = "(Hornem.) P.Karst."
$ws.Range("Q6").Value = 752290.4346396544
$ws.Range("R6").Value = 7092949.368184029

# Row 7
$ws.Range("A7").Value = 112027301
$ws.Range("B7").Value = 90658
$ws.Range("D7").Value = "NT"
$ws.Range("E7").Value = 4361
$ws.Range("F7").Value = "Orange taggsvamp"
$ws.Range("G7").Value = "Hydnellum aurantiacum"
$ws.Range("H7").Value = "(Batsch:Fr.) P.Karst."
$ws.Range("Q7").Value = 752329.5007247855
$ws.Range("R7").Value = 7092948.539329411

# Row 8
$ws.Range("A8").Value = 112027356
$ws.Range("B8").Value = 90666
$ws.Range("D8").Value = "LC"
$ws.Range("E8").Value = 4364
$ws.Range("F8").Value = "Dropptaggsvamp"
$ws.Range("G8").Value = "Hydnellum ferrugineum"
$ws.Range("H8").Value = "(Fr.:Fr.) P. Karst."
$ws.Range("Q8").Value = 752324.3875685094
$ws.Range("R8").Value = 7092929.537413944
$ws.Range("AC8").ClearContents()

# Row 9
$ws.Range("A9").Value = 112027275
$ws.Range("B9").Value = 90652
$ws.Range("D9").Value = "NT"
$ws.Range("E9").Value = 3100
$ws.Range("F9").Value = "Talltaggsvamp"
$ws.Range("G9").Value = "Bankera fuligineoalba"
$ws.Range("H9").Value = "(Schmidt : Fr.) Pouzar"
$ws.Range("Q9").Value = 752296.0450195302
$ws.Range("R9").Value = 7092962.213766729

# Row 10
$ws.Range("A10").Value = 112027366
$ws.Range("B10").Value = 90660
$ws.Range("D10").Value = "NT"
$ws.Range("E10").Value = 4362
$ws.Range("F10").Value = "Blå taggsvamp"
$ws.Range("G10").Value = "Hydnellum caeruleum"
$ws.Range("H10").Value = "(Hornem.) P.Karst."
$ws.Range("Q10").Value = 752248.2058430372
$ws.Range("R10").Value = 7093192.473407456

# Row 11
$ws.Range("A11").Value = 112027288
$ws.Range("B11").Value = 90678
$ws.Range("D11").Value = "LC"
$ws.Range("E11").Value = 4366
$ws.Range("F11").Value = "Skarp dropptaggsvamp"
$ws.Range("G11").Value = "Hydnellum peckii"
$ws.Range("H11").Value = "Banker"
$ws.Range("Q11").Value = 752281.5439615413
$ws.Range("R11").Value = 7092999.105565066
$ws.Range("AC11").ClearContents()

# Row 12
$ws.Range("A12").Value = 112027290
$ws.Range("B12").Value = 90658
$ws.Range("D12").Value = "NT"
$ws.Range("E12").Value = 4361
$ws.Range("F12").Value = "Orange taggsvamp"
$ws.Range("G12").Value = "Hydnellum aurantiacum"
$ws.Range("H12").Value = "(Batsch:Fr.) P.Karst."
$ws.Range("Q12").Value = 752268.0515137122
$ws.Range("R12").Value = 7093028.99707507
$ws.Range("AC12").Value = "Rikligt"

# Row 13
$ws.Range("A13").Value = 112027285
$ws.Range("B13").Value = 90682
$ws.Range("D13").Value = "NT"
$ws.Range("E13").Value = 2059
$ws.Range("F13").Value = "Skrovlig taggsvamp"
$ws.Range("G13").Value = "Hydnellum scabrosum"
$ws.Range("H13").Value = "(Fr.) E.Larss., K.H.Larss. & Kõljalg"
$ws.Range("Q13").Value = 752324.380454565
$ws.Range("R13").Value = 7092940.601741337
$ws.Range("AC13").ClearContents()

# Row 14
$ws.Range("A14").Value = 112027355
$ws.Range("B14").Value = 90666
$ws.Range("D14").Value = "LC"
$ws.Range("E14").Value = 4364
$ws.Range("F14").Value = "Dropptaggsvamp"
$ws.Range("G14").Value = "Hydnellum ferrugineum"
$ws.Range("H14").Value = "(Fr.:Fr.) P. Karst."
$ws.Range("Q14").Value = 752324.380454565
$ws.Range("R14").Value = 7092940.601741337

# Row 15
$ws.Range("A15").Value = 112027354
$ws.Range("B15").Value = 90666
$ws.Range("D15").Value = "LC"
$ws.Range("E15").Value = 4364
$ws.Range("F15").Value = "Dropptaggsvamp"
$ws.Range("G15").Value = "Hydnellum ferrugineum"
$ws.Range("H15").Value = "(Fr.:Fr.) P. Karst."
$ws.Range("Q15").Value = 752295.4558510378
$ws.Range("R15").Value = 7093002.441773332

# Row 16
$ws.Range("A16").Value = 112027360
$ws.Range("B16").Value = 90660
$ws.Range("D16").Value = "NT"
$ws.Range("E16").Value = 4362
$ws.Range("F16").Value = "Blå taggsvamp"
$ws.Range("G16").Value = "Hydnellum caeruleum"
$ws.Range("H16").Value = "(Hornem.) P.Karst."
$ws.Range("Q16").Value = 752296.5637017922
$ws.Range("R16").Value = 7093005.186741289
$ws.Range("AC16").Value = "ca 10 fruktkroppar"
